# DoInventory: rebuild the "Inventory - List of Articles" sheet with
# separate Title / Single Price / Quantity / Total Price columns instead
# of the old combined "Anzahl" / "Beschreibung" layout, and add a
# GRAND TOTAL row summing the total prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old A-column quantity values (data now lives in B:E).
$ws.Range("A1").ClearContents()
$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()

# Header row.
$ws.Range("B1").Value = "Title"
$ws.Range("C1").Value = "Single Price"
$ws.Range("D1").Value = "Quantity"
$ws.Range("E1").Value = "Total Price"

# RedBull line.
$ws.Range("B2").Value = "RedBull"
$ws.Range("C2").Value = "CHF 3.5"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 3.5

# Zweifel Chips line.
$ws.Range("B3").Value = "Zweifel Chips"
$ws.Range("C3").Value = "CHF 6.7"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 6.7

# Grand total row.
$ws.Range("A6").Value = "GRAND TOTAL:"
$ws.Range("E6").Value = 10.2
